# Updated Traceability Matrix with Use Case specification document names.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use-case codes (column B) and use-case names (column A) for rows 5..15
$codes = @("UC01","UC02","UC03","UC04","UC05","UC06","UC07","UC08","UC09","UC10","UC11")
$names = @(
    "Use Case - Modify Inventory",
    "Use Case - View Sales",
    "Use Case - Modify Route",
    "Use Case - Modify Truck",
    "Use Case - Process Batch File",
    "Use Case - Modify Item",
    "Use Case - Modify Driver",
    "Use Case - Modify Settings",
    "Use Case - Modify Voting",
    "Use Case - Modify Presets",
    "Use Case - View Fuel Usage"
)

for ($i = 0; $i -lt $codes.Length; $i++) {
    $row = 5 + $i
    $ws.Cells.Item($row, 2).Value = $codes[$i]
}

# Column A names were entered in this row order (matches the shared-string
# insertion order recorded by the original author's edit).
$nameRows = @(10,5,6,7,8,9,11,12,13,14,15)
foreach ($row in $nameRows) {
    $ws.Cells.Item($row, 1).Value = $names[$row - 5]
}

# Diagonal highlight additions: D2 and C3 pick up the gray-fill format used
# elsewhere on the diagonal (copy format from the existing diagonal cell C2).
$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("C2").Copy()
$ws.Range("C3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Update selection to match the new active cell
$ws.Range("F20").Select()
